$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the WHODAS_P_* variable names to WHODAS_SR_*.
# Rows 3-13 keep a 1:1 rename (01-11), row 14 previously duplicated
# "WHODAS_P_11" (a pre-existing data bug) and becomes the new
# "WHODAS_SR_12" entry, and the trailing "Days" variables shift down
# to match (Days01-03).
$ws.Range("B3").Value  = "WHODAS_SR_01"
$ws.Range("B4").Value  = "WHODAS_SR_02"
$ws.Range("B5").Value  = "WHODAS_SR_03"
$ws.Range("B6").Value  = "WHODAS_SR_04"
$ws.Range("B7").Value  = "WHODAS_SR_05"
$ws.Range("B8").Value  = "WHODAS_SR_06"
$ws.Range("B9").Value  = "WHODAS_SR_07"
$ws.Range("B10").Value = "WHODAS_SR_08"
$ws.Range("B11").Value = "WHODAS_SR_09"
$ws.Range("B12").Value = "WHODAS_SR_10"
$ws.Range("B13").Value = "WHODAS_SR_11"
$ws.Range("B14").Value = "WHODAS_SR_12"
$ws.Range("B15").Value = "WHODAS_SR_Days01"
$ws.Range("B16").Value = "WHODAS_SR_Days02"
$ws.Range("B17").Value = "WHODAS_SR_Days03"

# Reflect the updated selection / scroll position from the author's edit.
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B14").Select() | Out-Null
